$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial (45177 -> 45178) for every
# data row (rows 2 through 265). Bump it by one day for all of them.
for ($r = 2; $r -le 265; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
